$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4").Value = "2b"
$ws.Range("A6").Value = "2b"

$ws.Range("J5").Select()
